# Added to precis of Burke's Reflections
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New quote text added to the precis.
$newQuote = "if it be a thing confessed, that of such questions they cannot determine without rashness, inasmuch as a great part of them consisteth in special circumstances, and for one kind as many reasons may be brought as for another;"

# Row 6: add a 3rd "done" date (E6) matching existing C6/D6 date formatting.
$ws.Range("D6").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("E6").Value = "8/15/2018"

# Row 7: same as row 6.
$ws.Range("D7").Copy()
$ws.Range("E7").PasteSpecial(-4122)
$ws.Range("E7").Value = "8/15/2018"

# Row 8: add a 2nd date (D8).
$ws.Range("C8").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("D8").Value = "8/15/2018"

# Row 9: previously blank placeholder row -> now holds the new quote entry.
$ws.Range("B8:C8").Copy()
$ws.Range("B9:C9").PasteSpecial(-4122)
$ws.Range("B9").Value = $newQuote
$ws.Range("C9").Value = "8/15/2018"
$ws.Rows.Item(9).RowHeight = 60

# Move the active selection to C10 (next empty row), matching the author's
# cursor position after the edit.
$ws.Range("C10").Select()
